$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.7004376666666666
$ws.Range("H2").Value = 2.101313
$ws.Range("I2").Value = 0.04511966030063898
$ws.Range("J2").Value = 0.04511966030063898
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 5.900730666666667
$ws.Range("N2").Value = 17.702192
$ws.Range("O2").Value = 0.03970749001357476
$ws.Range("P2").Value = 0.03970749001357476
$ws.Range("Q2").Value = 4.133094019788444
$ws.Range("R2").Value = 37.19784617809599
$ws.Range("S2").Value = 0.001791588460803508
$ws.Range("T2").Value = 0.001791588460803508
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.7004376666666666
$ws.Range("H3").Value = 2.101313
$ws.Range("I3").Value = 0.04511966030063898
$ws.Range("J3").Value = 0.04511966030063898
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 111.5917106666667
$ws.Range("N3").Value = 334.775132
$ws.Range("O3").Value = 0.7509284844884279
$ws.Range("P3").Value = 0.7509284844884279
$ws.Range("Q3").Value = 78.16303743870176
$ws.Range("R3").Value = 703.4673369483158
$ws.Range("S3").Value = 0.03388163813019152
$ws.Range("T3").Value = 0.03388163813019151
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.7004376666666666
$ws.Range("H4").Value = 2.101313
$ws.Range("I4").Value = 0.04511966030063898
$ws.Range("J4").Value = 0.04511966030063898
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 31.11253633333333
$ws.Range("N4").Value = 93.337609
$ws.Range("O4").Value = 0.2093640254979974
$ws.Range("P4").Value = 0.2093640254979974
$ws.Range("Q4").Value = 21.79239235340189
$ws.Range("R4").Value = 196.131531180617
$ws.Range("S4").Value = 0.009446433709643962
$ws.Range("T4").Value = 0.009446433709643961
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.237305
$ws.Range("H5").Value = 0.711915
$ws.Range("I5").Value = 0.01528632952964618
$ws.Range("J5").Value = 0.01528632952964618
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.900730666666667
$ws.Range("N5").Value = 17.702192
$ws.Range("O5").Value = 0.03970749001357476
$ws.Range("P5").Value = 0.03970749001357476
$ws.Range("Q5").Value = 1.400272890853333
$ws.Range("R5").Value = 12.60245601768
$ws.Range("S5").Value = 0.0006069817771426387
$ws.Range("T5").Value = 0.0006069817771426386
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.237305
$ws.Range("H6").Value = 0.711915
$ws.Range("I6").Value = 0.01528632952964618
$ws.Range("J6").Value = 0.01528632952964618
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 111.5917106666667
$ws.Range("N6").Value = 334.775132
$ws.Range("O6").Value = 0.7509284844884279
$ws.Range("P6").Value = 0.7509284844884279
$ws.Range("Q6").Value = 26.48127089975333
$ws.Range("R6").Value = 238.33143809778
$ws.Range("S6").Value = 0.01147894026708791
$ws.Range("T6").Value = 0.01147894026708791
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.237305
$ws.Range("H7").Value = 0.711915
$ws.Range("I7").Value = 0.01528632952964618
$ws.Range("J7").Value = 0.01528632952964618
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 31.11253633333333
$ws.Range("N7").Value = 93.337609
$ws.Range("O7").Value = 0.2093640254979974
$ws.Range("P7").Value = 0.2093640254979974
$ws.Range("Q7").Value = 7.383160434581667
$ws.Range("R7").Value = 66.448443911235
$ws.Range("S7").Value = 0.003200407485415634
$ws.Range("T7").Value = 0.003200407485415633
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 14.58625866666667
$ws.Range("H8").Value = 43.758776
$ws.Range("I8").Value = 0.9395940101697148
$ws.Range("J8").Value = 0.9395940101697148
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 5.900730666666667
$ws.Range("N8").Value = 17.702192
$ws.Range("O8").Value = 0.03970749001357476
$ws.Range("P8").Value = 0.03970749001357476
$ws.Range("Q8").Value = 86.06958382633245
$ws.Range("R8").Value = 774.626254436992
$ws.Range("S8").Value = 0.03730891977562861
$ws.Range("T8").Value = 0.03730891977562861
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 14.58625866666667
$ws.Range("H9").Value = 43.758776
$ws.Range("I9").Value = 0.9395940101697148
$ws.Range("J9").Value = 0.9395940101697148
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 111.5917106666667
$ws.Range("N9").Value = 334.775132
$ws.Range("O9").Value = 0.7509284844884279
$ws.Range("P9").Value = 0.7509284844884279
$ws.Range("Q9").Value = 1627.705556839825
$ws.Range("R9").Value = 14649.35001155843
$ws.Range("S9").Value = 0.7055679060911484
$ws.Range("T9").Value = 0.7055679060911484
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 14.58625866666667
$ws.Range("H10").Value = 43.758776
$ws.Range("I10").Value = 0.9395940101697148
$ws.Range("J10").Value = 0.9395940101697148
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 31.11253633333333
$ws.Range("N10").Value = 93.337609
$ws.Range("O10").Value = 0.2093640254979974
$ws.Range("P10").Value = 0.2093640254979974
$ws.Range("Q10").Value = 453.8155027340649
$ws.Range("R10").Value = 4084.339524606584
$ws.Range("S10").Value = 0.1967171843029378
$ws.Range("T10").Value = 0.1967171843029378
